$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "95-29=66"
$t.Cell(1, 2).Range.Text = "42-27=15"
$t.Cell(1, 3).Range.Text = "31-5=26"
$t.Cell(1, 4).Range.Text = "91-2=89"
$t.Cell(1, 5).Range.Text = "25+62=87"
$t.Cell(2, 1).Range.Text = "51+35=86"
$t.Cell(2, 2).Range.Text = "14+9=23"
$t.Cell(2, 3).Range.Text = "5+82=87"
$t.Cell(2, 4).Range.Text = "72-6=66"
$t.Cell(2, 5).Range.Text = "5+33=38"
$t.Cell(3, 1).Range.Text = "95-31=64"
$t.Cell(3, 2).Range.Text = "99-85=14"
$t.Cell(3, 3).Range.Text = "60-20=40"
$t.Cell(3, 4).Range.Text = "49-0=49"
$t.Cell(3, 5).Range.Text = "21+53=74"
$t.Cell(4, 1).Range.Text = "45+9=54"
$t.Cell(4, 2).Range.Text = "52+26=78"
$t.Cell(4, 3).Range.Text = "37-11=26"
$t.Cell(4, 4).Range.Text = "19-18=1"
$t.Cell(4, 5).Range.Text = "16+20=36"
$t.Cell(5, 1).Range.Text = "7+39=46"
$t.Cell(5, 2).Range.Text = "71+20=91"
$t.Cell(5, 3).Range.Text = "86-69=17"
$t.Cell(5, 4).Range.Text = "94-74=20"
$t.Cell(5, 5).Range.Text = "11+44=55"
$t.Cell(6, 1).Range.Text = "34+57=91"
$t.Cell(6, 2).Range.Text = "62-39=23"
$t.Cell(6, 3).Range.Text = "50+31=81"
$t.Cell(6, 4).Range.Text = "78-72=6"
$t.Cell(6, 5).Range.Text = "56-23=33"
$t.Cell(7, 1).Range.Text = "96-12=84"
$t.Cell(7, 2).Range.Text = "52-41=11"
$t.Cell(7, 3).Range.Text = "9+13=22"
$t.Cell(7, 4).Range.Text = "0+55=55"
$t.Cell(7, 5).Range.Text = "81-15=66"
$t.Cell(8, 1).Range.Text = "45-30=15"
$t.Cell(8, 2).Range.Text = "31+17=48"
$t.Cell(8, 3).Range.Text = "38-34=4"
$t.Cell(8, 4).Range.Text = "85-9=76"
$t.Cell(8, 5).Range.Text = "72+2=74"
$t.Cell(9, 1).Range.Text = "88-64=24"
$t.Cell(9, 2).Range.Text = "25+6=31"
$t.Cell(9, 3).Range.Text = "13-13=0"
$t.Cell(9, 4).Range.Text = "73+25=98"
$t.Cell(9, 5).Range.Text = "26+28=54"
$t.Cell(10, 1).Range.Text = "24+71=95"
$t.Cell(10, 2).Range.Text = "6-0=6"
$t.Cell(10, 3).Range.Text = "96-80=16"
$t.Cell(10, 4).Range.Text = "72+13=85"
$t.Cell(10, 5).Range.Text = "31+50=81"
$t.Cell(11, 1).Range.Text = "43+37=80"
$t.Cell(11, 2).Range.Text = "4+14=18"
$t.Cell(11, 3).Range.Text = "49-12=37"
$t.Cell(11, 4).Range.Text = "32+44=76"
$t.Cell(11, 5).Range.Text = "34+47=81"
$t.Cell(12, 1).Range.Text = "73-47=26"
$t.Cell(12, 2).Range.Text = "34-1=33"
$t.Cell(12, 3).Range.Text = "18+78=96"
$t.Cell(12, 4).Range.Text = "79-4=75"
$t.Cell(12, 5).Range.Text = "7+80=87"
$t.Cell(13, 1).Range.Text = "0+57=57"
$t.Cell(13, 2).Range.Text = "60+33=93"
$t.Cell(13, 3).Range.Text = "26+25=51"
$t.Cell(13, 4).Range.Text = "72-57=15"
$t.Cell(13, 5).Range.Text = "29+42=71"
$t.Cell(14, 1).Range.Text = "33-21=12"
$t.Cell(14, 2).Range.Text = "27-19=8"
$t.Cell(14, 3).Range.Text = "0+47=47"
$t.Cell(14, 4).Range.Text = "83-45=38"
$t.Cell(14, 5).Range.Text = "92+2=94"
$t.Cell(15, 1).Range.Text = "39+48=87"
$t.Cell(15, 2).Range.Text = "82-9=73"
$t.Cell(15, 3).Range.Text = "66-35=31"
$t.Cell(15, 4).Range.Text = "46+18=64"
$t.Cell(15, 5).Range.Text = "64+19=83"
$t.Cell(16, 1).Range.Text = "62-16=46"
$t.Cell(16, 2).Range.Text = "87-45=42"
$t.Cell(16, 3).Range.Text = "42+8=50"
$t.Cell(16, 4).Range.Text = "51+2=53"
$t.Cell(16, 5).Range.Text = "54-1=53"
$t.Cell(17, 1).Range.Text = "14+20=34"
$t.Cell(17, 2).Range.Text = "61+10=71"
$t.Cell(17, 3).Range.Text = "50-33=17"
$t.Cell(17, 4).Range.Text = "32+44=76"
$t.Cell(17, 5).Range.Text = "40+8=48"
$t.Cell(18, 1).Range.Text = "5+39=44"
$t.Cell(18, 2).Range.Text = "66-23=43"
$t.Cell(18, 3).Range.Text = "55-31=24"
$t.Cell(18, 4).Range.Text = "28-7=21"
$t.Cell(18, 5).Range.Text = "43-30=13"
$t.Cell(19, 1).Range.Text = "84-36=48"
$t.Cell(19, 2).Range.Text = "75-22=53"
$t.Cell(19, 3).Range.Text = "68-10=58"
$t.Cell(19, 4).Range.Text = "78-46=32"
$t.Cell(19, 5).Range.Text = "30+26=56"
$t.Cell(20, 1).Range.Text = "45-15=30"
$t.Cell(20, 2).Range.Text = "58+22=80"
$t.Cell(20, 3).Range.Text = "15+0=15"
$t.Cell(20, 4).Range.Text = "7+0=7"
$t.Cell(20, 5).Range.Text = "65-61=4"
